# Update patient sheet: row 3 values change, row 4 gets fully cleared
$wsPatient = $excel.ActiveWorkbook.Worksheets.Item("patient")
$wsPatient.Range("A3").Value = "CRC0228PRaS"
$wsPatient.Range("B3").Value = "not provided"
$wsPatient.Range("A4:B4").Clear()

# Update pdx_model sheet: model_id values corrected
$wsPdx = $excel.ActiveWorkbook.Worksheets.Item("pdx_model")
$wsPdx.Range("A2").Value = "CRC0228PR"
$wsPdx.Range("A3").Value = "CRC0228PRaS"

# Update model_validation sheet: model_id value corrected
$wsModelValidation = $excel.ActiveWorkbook.Worksheets.Item("model_validation")
$wsModelValidation.Range("A3").Value = "CRC0228PRaS"

# Update cell_model sheet: model_id value made unique
$wsCellModel = $excel.ActiveWorkbook.Worksheets.Item("cell_model")
$wsCellModel.Range("A3").Value = "CRC0014LM_2"

# Restore selections on each sheet to match the saved workbook state
$wsPdx.Range("A3").Select()
$wsModelValidation.Range("A3").Select()
$wsCellModel.Range("A3").Select()

$wsPatientSample = $excel.ActiveWorkbook.Worksheets.Item("patient_sample")
$wsPatientSample.Range("T3").Select()

$wsSharing = $excel.ActiveWorkbook.Worksheets.Item("sharing")
$wsSharing.Range("D3").Select()

# Finish on the patient sheet with the B4 selection, matching the saved file
$wsPatient.Range("B4").Select()
